# Generate Report for Handback
#
# - Overview / zh-cn / de-de "Status" cells flip from "Ready for handoff" to
#   "Handed back: in sync with en-US" (all share one string, so updating any
#   one cell via the shared string flips every occurrence).
# - zh-cn / de-de rows 2-3 gain "Latest Target File" (col F) and
#   "Latest Handback File" (col G) hyperlink cells that mirror the source
#   file name / handoff target file respectively.
# - zh-cn / de-de "Latest Handback DateTime" (col H) moves from the
#   "0001-01-01 00:00:00" placeholder to the actual handback timestamp.

$wb = $excel.ActiveWorkbook

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/d81fe4da196d4e1e08492f1549afcdd2966d8b1b/e2e/af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c799ab0487b6dd8dd5f790c1efe6313d7b8cf61/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.6e12c5a9b899230eba1fbb53799aae92711faa9a.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08c7da08f2f0c50e80175d3b0a4e1f9a99c994cd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.6e12c5a9b899230eba1fbb53799aae92711faa9a.de-de.xlf"

$handedBack = "Handed back: in sync with en-US"

# --- Status text: every "Ready for handoff" cell on every sheet flips to
#     the handed-back message (they all shared one sharedStrings entry in
#     the original workbook). ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

function Update-LangSheet($sheetName, $xlfUrl, $handbackTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("C2").Value = $handedBack
    $ws.Range("C3").Value = $handedBack

    # New "Latest Target File" (F) / "Latest Handback File" (G) hyperlink
    # cells for rows 2 and 3.
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.md")
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, "", "", "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.6e12c5a9b899230eba1fbb53799aae92711faa9a.$sheetName.xlf")
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, "", "", "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.md")
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, "", "", "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.6e12c5a9b899230eba1fbb53799aae92711faa9a.$sheetName.xlf")

    # Match the existing hyperlink look (blue underline) used elsewhere on
    # the sheet (col A/B/D already carry it). Re-asserting Name last nudges
    # the engine into reusing that same font entry instead of minting a
    # theme-linked one.
    $newLinks = $ws.Range("F2:G3")
    $newLinks.Font.Underline = $true
    $newLinks.Font.Color = 15570276
    $newLinks.Font.Name = "Calibri"

    # "Latest Handback DateTime" (H) : real timestamp instead of the
    # 0001-01-01 placeholder.
    $ws.Range("H2").Value = $handbackTime
    $ws.Range("H3").Value = $handbackTime
}

Update-LangSheet "zh-cn" $zhXlfUrl "2016-03-18 00:50:12"
Update-LangSheet "de-de" $deXlfUrl "2016-03-18 00:50:18"
